$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NAF")

function Set-TextValue($range, $value) {
    # Forces a cell to hold a plain string even when the text looks
    # numeric (e.g. "59.9241139"), without leaving a stray NumberFormat
    # behind once we're done (mirrors the General/Normal style that the
    # surrounding cells already use).
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

# --- New row 9: NAF Senter Oslo -------------------------------------------
$ws.Range("A9").Value2 = "NAF Senter Oslo"
$ws.Range("B9").Value2 = "Eikenga 9"
$ws.Range("C9").Value2 = "Oslo"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0579"
Set-TextValue $ws.Range("E9") "59.9241139"
Set-TextValue $ws.Range("F9") "10.8047797"

# --- New row 10: NAF Senter Rykkinn ----------------------------------------
$ws.Range("A10").Value2 = "NAF Senter Rykkinn"
$ws.Range("B10").Value2 = "Folkvangveien 22"
$ws.Range("C10").Value2 = "Rykkinn"
$ws.Range("D10").Value2 = 1348
$ws.Range("D10").NumberFormat = "@"
Set-TextValue $ws.Range("E10") "59.9244041442871"
Set-TextValue $ws.Range("F10") "10.4956493377686"

# --- Rename the existing NAF stations in column A --------------------------
$ws.Range("A8").Value2 = "NAF Senter Kistiansand"
$ws.Range("A7").Value2 = "NAF Senter Steinkjær"
$ws.Range("A6").Value2 = "NAF Senter Arendal"
$ws.Range("A5").Value2 = "NAF Senter Namsos"
$ws.Range("A4").Value2 = "NAF Senter Bergen"
$ws.Range("A2").Value2 = "NAF Senter Trondheim"
$ws.Range("A3").Value2 = "NAF Senter Sandnes"

# --- Widen column A now that the station names are longer ------------------
# ColumnWidth and the stored OOXML character width differ by a constant
# offset (~5/6 of a character) in this engine, so back that out to land on
# an on-disk width of exactly 20.
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668

# --- Make NAF the active sheet, with D10 selected, like the saved file ----
$ws.Activate() | Out-Null
$ws.Range("D10").Select() | Out-Null

Write-Host "NAF sheet updated with Oslo/Rykkinn stations and renamed centers"
